$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) to be treated as text so values like "1.001" or
# "0.00001047" are not reinterpreted as numbers and lose their exact
# textual representation (trailing zeros, decimal grouping, etc.).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.971.46"
$ws.Range("E2").Value = "  +2.27%  "
$ws.Range("D3").Value = "1.676.57"
$ws.Range("E3").Value = "  +1.40%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "328.87"
$ws.Range("E5").Value = "  +6.97%  "
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "0.3658"
$ws.Range("E7").Value = "  +1.09%  "
$ws.Range("D8").Value = "46.63"
$ws.Range("E8").Value = "  -1.75%  "
$ws.Range("D9").Value = "0.3245"
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("E10").Value = "  +1.44%  "
$ws.Range("D11").Value = "0.07075"
$ws.Range("E11").Value = "  +1.76%  "
$ws.Range("D12").Value = "0.9999"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "6.079"
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("D14").Value = "19.62"
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("D15").Value = "1.676.49"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").Value = "6.626"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "0.00001047"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "0.06582"
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("D19").Value = "0.9999"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "78.94"
$ws.Range("E20").Value = "  +2.94%  "
$ws.Range("D21").Value = "15.90"
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "12.94"
$ws.Range("E23").Value = "  +2.04%  "
$ws.Range("D24").Value = "24.964.62"
$ws.Range("E24").Value = "  +2.25%  "
$ws.Range("D25").Value = "2.440"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").Value = "2.400"
$ws.Range("E26").Value = "  +2.03%  "
$ws.Range("D27").Value = "148.17"
$ws.Range("E27").Value = "  +1.29%  "
$ws.Range("D28").Value = "18.76"
$ws.Range("E28").Value = "  +2.02%  "
$ws.Range("D29").Value = "1.863.74"
$ws.Range("E29").Value = "  +1.43%  "
$ws.Range("D30").Value = "125.50"
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("D31").Value = "1.185"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").Value = "4.074"
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("D33").Value = "5.780"
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("D34").Value = "0.08480"
$ws.Range("E34").Value = "  +1.80%  "
$ws.Range("D35").Value = "1.651"
$ws.Range("E35").Value = "  -0.97%  "
$ws.Range("D36").Value = "12.30"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "5.176"
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("D38").Value = "0.02253"
$ws.Range("E38").Value = "  +2.04%  "
$ws.Range("E39").Value = "  +2.02%  "
$ws.Range("D40").Value = "0.06027"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("D41").Value = "0.2094"
$ws.Range("E41").Value = "  +2.08%  "
$ws.Range("D42").Value = "8.241"
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("D43").Value = "0.9989"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").Value = "0.5963"
$ws.Range("E44").Value = "  +2.06%  "
$ws.Range("D45").Value = "13.76"
$ws.Range("E45").Value = "  +9.51%  "
$ws.Range("D46").Value = "3.841"
$ws.Range("E46").Value = "  +3.03%  "
$ws.Range("D47").Value = "0.5756"
$ws.Range("E47").Value = "  +3.24%  "
$ws.Range("D48").Value = "125.44"
$ws.Range("E48").Value = "  +2.76%  "
$ws.Range("D49").Value = "1.967"
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("D50").Value = "0.07019"
$ws.Range("E50").Value = "  +1.57%  "
$ws.Range("E51").Value = "  +3.02%  "

# Restore the original (default) cell style now that the text values are set,
# so no visible/style attribute changes leak into the cells.
$ws.Range("D2:D51").Style = "Normal"
